$wb = $excel.ActiveWorkbook

# --- DeveloperTabData sheet: insert a new top row for a new test run ---
$wsDev = $wb.Worksheets.Item("DeveloperTabData")

# Insert a new row 1; existing rows 1-6 shift down to 2-7.
$wsDev.Rows.Item(1).Insert()

# Copy number formatting (style) from the (now shifted) row below so the
# new date cells pick up the same date/time style instead of a generic one.
$wsDev.Range("E2:F2").Copy()
$wsDev.Range("E1:F1").PasteSpecial(-4122)

# Populate the new row's values.
$wsDev.Cells.Item(1, 1).Value = 321
$wsDev.Cells.Item(1, 2).Value = "testautocomponent_379767"
$wsDev.Cells.Item(1, 3).Value = "Shell Script"
$wsDev.Cells.Item(1, 4).Value = "Suyog Talathi"
$wsDev.Cells.Item(1, 5).Value = 43210.471770833334
$wsDev.Cells.Item(1, 7).Value = "Created"

# --- Make "Visualize" the active/selected sheet instead of "DeveloperTabData" ---
$wsVisualize = $wb.Worksheets.Item("Visualize")
$wsVisualize.Activate()
